$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dll1"
$ws.Cells.Item(2, 3).Value = "Notch2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 16.39002933333333
$ws.Cells.Item(2, 8).Value = 49.170088
$ws.Cells.Item(2, 9).Value = 0.5551882184054378
$ws.Cells.Item(2, 10).Value = 0.5551882184054378
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.174933333333333
$ws.Cells.Item(2, 14).Value = 3.5248
$ws.Cells.Item(2, 15).Value = 0.01171850713626266
$ws.Cells.Item(2, 16).Value = 0.01171850713626266
$ws.Cells.Item(2, 17).Value = 19.25719179804445
$ws.Cells.Item(2, 18).Value = 173.3147261824
$ws.Cells.Item(2, 19).Value = 0.006505977099353073
$ws.Cells.Item(2, 20).Value = 0.006505977099353072

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Dll1"
$ws.Cells.Item(3, 3).Value = "Notch2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 16.39002933333333
$ws.Cells.Item(3, 8).Value = 49.170088
$ws.Cells.Item(3, 9).Value = 0.5551882184054378
$ws.Cells.Item(3, 10).Value = 0.5551882184054378
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 27.50472933333333
$ws.Cells.Item(3, 14).Value = 82.514188
$ws.Cells.Item(3, 15).Value = 0.2743256641287217
$ws.Cells.Item(3, 16).Value = 0.2743256641287218
$ws.Cells.Item(3, 17).Value = 450.8033205787272
$ws.Cells.Item(3, 18).Value = 4057.229885208544
$ws.Cells.Item(3, 19).Value = 0.1523023767305135
$ws.Cells.Item(3, 20).Value = 0.1523023767305136

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Dll1"
$ws.Cells.Item(4, 3).Value = "Notch2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 16.39002933333333
$ws.Cells.Item(4, 8).Value = 49.170088
$ws.Cells.Item(4, 9).Value = 0.5551882184054378
$ws.Cells.Item(4, 10).Value = 0.5551882184054378
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 39.361408
$ws.Cells.Item(4, 14).Value = 118.084224
$ws.Cells.Item(4, 15).Value = 0.3925813724534833
$ws.Cells.Item(4, 16).Value = 0.3925813724534833
$ws.Cells.Item(4, 17).Value = 645.1346317213015
$ws.Cells.Item(4, 18).Value = 5806.211685491712
$ws.Cells.Item(4, 19).Value = 0.217956552751611
$ws.Cells.Item(4, 20).Value = 0.217956552751611

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Dll1"
$ws.Cells.Item(5, 3).Value = "Notch2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 16.39002933333333
$ws.Cells.Item(5, 8).Value = 49.170088
$ws.Cells.Item(5, 9).Value = 0.5551882184054378
$ws.Cells.Item(5, 10).Value = 0.5551882184054378
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 32.221985
$ws.Cells.Item(5, 14).Value = 96.665955
$ws.Cells.Item(5, 15).Value = 0.3213744562815322
$ws.Cells.Item(5, 16).Value = 0.3213744562815322
$ws.Cells.Item(5, 17).Value = 528.1192793282266
$ws.Cells.Item(5, 18).Value = 4753.07351395404
$ws.Cells.Item(5, 19).Value = 0.1784233118239601
$ws.Cells.Item(5, 20).Value = 0.1784233118239601

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Dll1"
$ws.Cells.Item(6, 3).Value = "Notch2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.2120556666666667
$ws.Cells.Item(6, 8).Value = 0.636167
$ws.Cells.Item(6, 9).Value = 0.007183074867352934
$ws.Cells.Item(6, 10).Value = 0.007183074867352935
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.174933333333333
$ws.Cells.Item(6, 14).Value = 3.5248
$ws.Cells.Item(6, 15).Value = 0.01171850713626266
$ws.Cells.Item(6, 16).Value = 0.01171850713626266
$ws.Cells.Item(6, 17).Value = 0.2491512712888889
$ws.Cells.Item(6, 18).Value = 2.2423614416
$ws.Cells.Item(6, 19).Value = 0.00008417491409338429
$ws.Cells.Item(6, 20).Value = 0.0000841749140933843

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Dll1"
$ws.Cells.Item(7, 3).Value = "Notch2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.2120556666666667
$ws.Cells.Item(7, 8).Value = 0.636167
$ws.Cells.Item(7, 9).Value = 0.007183074867352934
$ws.Cells.Item(7, 10).Value = 0.007183074867352935
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 27.50472933333333
$ws.Cells.Item(7, 14).Value = 82.514188
$ws.Cells.Item(7, 15).Value = 0.2743256641287217
$ws.Cells.Item(7, 16).Value = 0.2743256641287218
$ws.Cells.Item(7, 17).Value = 5.832533715266222
$ws.Cells.Item(7, 18).Value = 52.49280343739601
$ws.Cells.Item(7, 19).Value = 0.001970501783472923
$ws.Cells.Item(7, 20).Value = 0.001970501783472924

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Dll1"
$ws.Cells.Item(8, 3).Value = "Notch2"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.2120556666666667
$ws.Cells.Item(8, 8).Value = 0.636167
$ws.Cells.Item(8, 9).Value = 0.007183074867352934
$ws.Cells.Item(8, 10).Value = 0.007183074867352935
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 39.361408
$ws.Cells.Item(8, 14).Value = 118.084224
$ws.Cells.Item(8, 15).Value = 0.3925813724534833
$ws.Cells.Item(8, 16).Value = 0.3925813724534833
$ws.Cells.Item(8, 17).Value = 8.346809614378667
$ws.Cells.Item(8, 18).Value = 75.12128652940801
$ws.Cells.Item(8, 19).Value = 0.002819941389861537
$ws.Cells.Item(8, 20).Value = 0.002819941389861538

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Dll1"
$ws.Cells.Item(9, 3).Value = "Notch2"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.2120556666666667
$ws.Cells.Item(9, 8).Value = 0.636167
$ws.Cells.Item(9, 9).Value = 0.007183074867352934
$ws.Cells.Item(9, 10).Value = 0.007183074867352935
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 32.221985
$ws.Cells.Item(9, 14).Value = 96.665955
$ws.Cells.Item(9, 15).Value = 0.3213744562815322
$ws.Cells.Item(9, 16).Value = 0.3213744562815322
$ws.Cells.Item(9, 17).Value = 6.832854510498333
$ws.Cells.Item(9, 18).Value = 61.495690594485
$ws.Cells.Item(9, 19).Value = 0.002308456779925088
$ws.Cells.Item(9, 20).Value = 0.002308456779925089

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Dll1"
$ws.Cells.Item(10, 3).Value = "Notch2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 12.89091033333333
$ws.Cells.Item(10, 8).Value = 38.672731
$ws.Cells.Item(10, 9).Value = 0.4366606914505164
$ws.Cells.Item(10, 10).Value = 0.4366606914505165
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.174933333333333
$ws.Cells.Item(10, 14).Value = 3.5248
$ws.Cells.Item(10, 15).Value = 0.01171850713626266
$ws.Cells.Item(10, 16).Value = 0.01171850713626266
$ws.Cells.Item(10, 17).Value = 15.14596024764444
$ws.Cells.Item(10, 18).Value = 136.3136422288
$ws.Cells.Item(10, 19).Value = 0.005117011428888263
$ws.Cells.Item(10, 20).Value = 0.005117011428888263

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Dll1"
$ws.Cells.Item(11, 3).Value = "Notch2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 12.89091033333333
$ws.Cells.Item(11, 8).Value = 38.672731
$ws.Cells.Item(11, 9).Value = 0.4366606914505164
$ws.Cells.Item(11, 10).Value = 0.4366606914505165
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 27.50472933333333
$ws.Cells.Item(11, 14).Value = 82.514188
$ws.Cells.Item(11, 15).Value = 0.2743256641287217
$ws.Cells.Item(11, 16).Value = 0.2743256641287218
$ws.Cells.Item(11, 17).Value = 354.5609995786031
$ws.Cells.Item(11, 18).Value = 3191.048996207428
$ws.Cells.Item(11, 19).Value = 0.1197872341810697
$ws.Cells.Item(11, 20).Value = 0.1197872341810698

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Dll1"
$ws.Cells.Item(12, 3).Value = "Notch2"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 12.89091033333333
$ws.Cells.Item(12, 8).Value = 38.672731
$ws.Cells.Item(12, 9).Value = 0.4366606914505164
$ws.Cells.Item(12, 10).Value = 0.4366606914505165
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 39.361408
$ws.Cells.Item(12, 14).Value = 118.084224
$ws.Cells.Item(12, 15).Value = 0.3925813724534833
$ws.Cells.Item(12, 16).Value = 0.3925813724534833
$ws.Cells.Item(12, 17).Value = 507.4043811217493
$ws.Cells.Item(12, 18).Value = 4566.639430095745
$ws.Cells.Item(12, 19).Value = 0.1714248535461307
$ws.Cells.Item(12, 20).Value = 0.1714248535461308

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Dll1"
$ws.Cells.Item(13, 3).Value = "Notch2"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 12.89091033333333
$ws.Cells.Item(13, 8).Value = 38.672731
$ws.Cells.Item(13, 9).Value = 0.4366606914505164
$ws.Cells.Item(13, 10).Value = 0.4366606914505165
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 32.221985
$ws.Cells.Item(13, 14).Value = 96.665955
$ws.Cells.Item(13, 15).Value = 0.3213744562815322
$ws.Cells.Item(13, 16).Value = 0.3213744562815322
$ws.Cells.Item(13, 17).Value = 415.3707193970116
$ws.Cells.Item(13, 18).Value = 3738.336474573105
$ws.Cells.Item(13, 19).Value = 0.1403315922944276
$ws.Cells.Item(13, 20).Value = 0.1403315922944277

# Row 14
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Dll1"
$ws.Cells.Item(14, 3).Value = "Notch2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.02857733333333333
$ws.Cells.Item(14, 8).Value = 0.085732
$ws.Cells.Item(14, 9).Value = 0.0009680152766929151
$ws.Cells.Item(14, 10).Value = 0.0009680152766929153
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.174933333333333
$ws.Cells.Item(14, 14).Value = 3.5248
$ws.Cells.Item(14, 15).Value = 0.01171850713626266
$ws.Cells.Item(14, 16).Value = 0.01171850713626266
$ws.Cells.Item(14, 17).Value = 0.03357646151111111
$ws.Cells.Item(14, 18).Value = 0.3021881536
$ws.Cells.Item(14, 19).Value = 0.0000113436939279372
$ws.Cells.Item(14, 20).Value = 0.0000113436939279372

# Row 15
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Dll1"
$ws.Cells.Item(15, 3).Value = "Notch2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.02857733333333333
$ws.Cells.Item(15, 8).Value = 0.085732
$ws.Cells.Item(15, 9).Value = 0.0009680152766929151
$ws.Cells.Item(15, 10).Value = 0.0009680152766929153
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 27.50472933333333
$ws.Cells.Item(15, 14).Value = 82.514188
$ws.Cells.Item(15, 15).Value = 0.2743256641287217
$ws.Cells.Item(15, 16).Value = 0.2743256641287218
$ws.Cells.Item(15, 17).Value = 0.7860118184017778
$ws.Cells.Item(15, 18).Value = 7.074106365616
$ws.Cells.Item(15, 19).Value = 0.0002655514336655322
$ws.Cells.Item(15, 20).Value = 0.0002655514336655324

# Row 16
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Dll1"
$ws.Cells.Item(16, 3).Value = "Notch2"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.02857733333333333
$ws.Cells.Item(16, 8).Value = 0.085732
$ws.Cells.Item(16, 9).Value = 0.0009680152766929151
$ws.Cells.Item(16, 10).Value = 0.0009680152766929153
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 39.361408
$ws.Cells.Item(16, 14).Value = 118.084224
$ws.Cells.Item(16, 15).Value = 0.3925813724534833
$ws.Cells.Item(16, 16).Value = 0.3925813724534833
$ws.Cells.Item(16, 17).Value = 1.124844076885333
$ws.Cells.Item(16, 18).Value = 10.123596691968
$ws.Cells.Item(16, 19).Value = 0.0003800247658800431
$ws.Cells.Item(16, 20).Value = 0.0003800247658800432

# Row 17
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Dll1"
$ws.Cells.Item(17, 3).Value = "Notch2"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.02857733333333333
$ws.Cells.Item(17, 8).Value = 0.085732
$ws.Cells.Item(17, 9).Value = 0.0009680152766929151
$ws.Cells.Item(17, 10).Value = 0.0009680152766929153
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 32.221985
$ws.Cells.Item(17, 14).Value = 96.665955
$ws.Cells.Item(17, 15).Value = 0.3213744562815322
$ws.Cells.Item(17, 16).Value = 0.3213744562815322
$ws.Cells.Item(17, 17).Value = 0.9208184060066665
$ws.Cells.Item(17, 18).Value = 8.28736565406
$ws.Cells.Item(17, 19).Value = 0.0003110953832194025
$ws.Cells.Item(17, 20).Value = 0.0003110953832194026
